{"js": "// Update the worksheet date and every \"AxB=\" multiplication prompt to the\n// new values per the commit diff. Each old string is unique in the document,\n// so a simple search + replace (matchCase, no wildcards) is safe.\nconst replacements = [\n  [\"2023-09-06 Wednesday\", \"2023-09-07 Thursday\"],\n  [\"46\u00d767=\", \"79\u00d797=\"],\n  [\"69\u00d746=\", \"20\u00d754=\"],\n  [\"20\u00d794=\", \"33\u00d766=\"],\n  [\"37\u00d713=\", \"45\u00d715=\"],\n  [\"29\u00d760=\", \"84\u00d768=\"],\n  [\"86\u00d799=\", \"86\u00d757=\"],\n  [\"92\u00d757=\", \"54\u00d765=\"],\n  [\"16\u00d732=\", \"63\u00d797=\"],\n  [\"46\u00d752=\", \"16\u00d763=\"],\n  [\"55\u00d721=\", \"15\u00d749=\"],\n  [\"64\u00d781=\", \"70\u00d788=\"],\n  [\"21\u00d768=\", \"73\u00d774=\"],\n  [\"34\u00d723=\", \"27\u00d792=\"],\n  [\"60\u00d758=\", \"52\u00d775=\"],\n  [\"68\u00d756=\", \"49\u00d760=\"],\n  [\"72\u00d732=\", \"89\u00d736=\"],\n  [\"49\u00d729=\", \"19\u00d763=\"],\n  [\"28\u00d736=\", \"36\u00d767=\"],\n  [\"59\u00d731=\", \"31\u00d745=\"],\n  [\"68\u00d726=\", \"30\u00d767=\"],\n  [\"95\u00d759=\", \"28\u00d767=\"],\n  [\"51\u00d711=\", \"85\u00d745=\"],\n  [\"46\u00d797=\", \"66\u00d787=\"],\n  [\"65\u00d736=\", \"38\u00d785=\"],\n  [\"65\u00d785=\", \"28\u00d734=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"AxB=\" multiplication prompt to the\n# new values per the commit diff. Each old string is unique in the document,\n# so Find/Replace (MatchCase, MatchWholeWord off, Replace:=wdReplaceAll) is\n# safe for every entry.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2023-09-06 Wednesday\", \"2023-09-07 Thursday\"),\n  @(\"46\u00d767=\", \"79\u00d797=\"),\n  @(\"69\u00d746=\", \"20\u00d754=\"),\n  @(\"20\u00d794=\", \"33\u00d766=\"),\n  @(\"37\u00d713=\", \"45\u00d715=\"),\n  @(\"29\u00d760=\", \"84\u00d768=\"),\n  @(\"86\u00d799=\", \"86\u00d757=\"),\n  @(\"92\u00d757=\", \"54\u00d765=\"),\n  @(\"16\u00d732=\", \"63\u00d797=\"),\n  @(\"46\u00d752=\", \"16\u00d763=\"),\n  @(\"55\u00d721=\", \"15\u00d749=\"),\n  @(\"64\u00d781=\", \"70\u00d788=\"),\n  @(\"21\u00d768=\", \"73\u00d774=\"),\n  @(\"34\u00d723=\", \"27\u00d792=\"),\n  @(\"60\u00d758=\", \"52\u00d775=\"),\n  @(\"68\u00d756=\", \"49\u00d760=\"),\n  @(\"72\u00d732=\", \"89\u00d736=\"),\n  @(\"49\u00d729=\", \"19\u00d763=\"),\n  @(\"28\u00d736=\", \"36\u00d767=\"),\n  @(\"59\u00d731=\", \"31\u00d745=\"),\n  @(\"68\u00d726=\", \"30\u00d767=\"),\n  @(\"95\u00d759=\", \"28\u00d767=\"),\n  @(\"51\u00d711=\", \"85\u00d745=\"),\n  @(\"46\u00d797=\", \"66\u00d787=\"),\n  @(\"65\u00d736=\", \"38\u00d785=\"),\n  @(\"65\u00d785=\", \"28\u00d734=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $new\n\n  $found = $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    throw \"Not found: $old\"\n  }\n}\n"}
